# "Generate Report for Handback" — mark the two handed-off files as handed
# back: in sync with en-US, and record the generated target/handback files
# + timestamps for each locale sheet.

$wb = $excel.ActiveWorkbook

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2325406c731ef4c32eeda26c9f1e2708183c183/e2e/"
$file1 = "0f79903a-6429-4909-b05e-70595a75186a"
$file2 = "8405c3c3-f951-44bf-91aa-bf8146efbc86"
$file1Md = "$file1.md"
$file2Md = "$file2.md"
$file1Url = "$mdUrlBase$file1Md"
$file2Url = "$mdUrlBase$file2Md"

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: flip both rows' status columns to "handed back" ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# --- Per-locale detail sheets: zh-cn, de-de -----------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Suffix = "zh-cn"; K2 = "2016-11-02 05:41:41"; K3 = "2016-11-02 05:41:41" },
    @{ Sheet = "de-de"; Suffix = "de-de"; K2 = "2016-11-02 05:41:59"; K3 = "2016-11-02 05:41:59" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Status column (shared string with Overview E/F, updates in lockstep)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Row 2 — file1
    $ws.Hyperlinks.Add($ws.Range("I2"), $file1Url, [System.Type]::Missing, [System.Type]::Missing, $file1Md) | Out-Null
    $ws.Range("J2").Value = "$file1.e3f0eb7c0e947d8a66aec2f7bc08f24423ab0bfc.$($loc.Suffix).xlf"
    $ws.Range("K2").Value = $loc.K2

    # Row 3 — file2
    $ws.Hyperlinks.Add($ws.Range("I3"), $file2Url, [System.Type]::Missing, [System.Type]::Missing, $file2Md) | Out-Null
    $ws.Range("J3").Value = "$file2.a3dfacdfbaa86dced21afe30143214db96cfd395.$($loc.Suffix).xlf"
    $ws.Range("K3").Value = $loc.K3

    # Column widths: Status (C) widens like Overview's E/F; Latest Target
    # File (I) and Latest Handback File (J) both grow to fit full filenames.
    $ws.Columns.Item(3).ColumnWidth = 29.1666666666667
    $ws.Columns.Item(9).ColumnWidth = 39.1666666666667
    $ws.Columns.Item(10).ColumnWidth = 39.1666666666667
}

Write-Output "Generate Report for Handback: done"
